$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The price strings (e.g. "470.000€") look numeric/currency to Excel's
# auto-detection, which would otherwise coerce them into numeric cells.
# Temporarily force the cells to Text so the values are written back as
# literal strings (matching the original inlineStr cells), then clear the
# number-format override again so the cells end up with no special style,
# just like in the source file.
$rng = $ws.Range("C4:C7")
$rng.NumberFormat = "@"

$ws.Range("C4").Value = "470.000€"
$ws.Range("C5").Value = "395.000€"
$ws.Range("C6").Value = "490.000€"
$ws.Range("C7").Value = "430.000€"

$rng.ClearFormats()
